$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C (names/links) are plain text; set directly.
# Column D (price) and E (volume) must be forced to text so Excel
# does not reinterpret numeric-looking strings (e.g. '5.33') as numbers.

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("B50").Value = 'Polygon'
$ws.Range("C50").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'

$textCells = @{
    'D2' = '56.411.82'
    'E2' = '  +2.47%  '
    'D3' = '2.316.86'
    'E3' = '  +1.74%  '
    'E4' = '  +0.15%  '
    'D5' = '514.95'
    'E5' = '  +1.66%  '
    'D6' = '134.81'
    'E6' = '  +4.90%  '
    'D7' = '0.998'
    'E7' = '  +0.27%  '
    'D8' = '0.536'
    'E8' = '  +1.15%  '
    'D9' = '2.329.18'
    'E9' = '  +1.65%  '
    'E10' = '  +3.22%  '
    'E11' = '  -1.29%  '
    'D12' = '5.33'
    'E12' = '  +5.07%  '
    'D13' = '0.342'
    'E13' = '  +0.10%  '
    'D14' = '23.80'
    'E14' = '  +0.64%  '
    'D15' = '2.734.76'
    'E15' = '  +1.92%  '
    'D16' = '56.550.64'
    'E16' = '  +2.70%  '
    'D17' = '0.0000134'
    'E17' = '  +2.01%  '
    'D18' = '2.321.50'
    'E18' = '  +1.85%  '
    'D19' = '10.48'
    'E19' = '  +1.05%  '
    'D20' = '325.36'
    'E20' = '  +3.50%  '
    'D21' = '4.21'
    'E21' = '  +0.20%  '
    'D22' = '6.57'
    'E22' = '  +0.18%  '
    'E23' = '  +0.24%  '
    'D24' = '60.59'
    'E24' = '  +1.21%  '
    'E25' = '  +6.19%  '
    'D26' = '0.996'
    'E26' = '  +0.16%  '
    'D27' = '7.94'
    'E27' = '  +5.07%  '
    'D28' = '1.27'
    'E28' = '  +10.86%  '
    'D29' = '169.12'
    'E29' = '  -1.17%  '
    'D30' = '0.0₃0738'
    'E30' = '  +4.37%  '
    'E31' = '  +2.73%  '
    'D32' = '6.17'
    'E32' = '  +0.37%  '
    'D33' = '18.46'
    'E33' = '  +2.60%  '
    'E34' = '  +0.00%  '
    'D35' = '0.995'
    'E35' = '  +0.14%  '
    'D36' = '1.25'
    'E36' = '  +1.30%  '
    'D37' = '0.913'
    'E37' = '  +1.18%  '
    'D38' = '3.99'
    'E38' = '  +2.49%  '
    'D39' = '1.55'
    'E39' = '  +6.57%  '
    'D40' = '38.26'
    'E40' = '  +3.90%  '
    'D41' = '0.380'
    'E41' = '  +1.61%  '
    'D42' = '141.14'
    'E42' = '  +3.22%  '
    'D43' = '3.59'
    'E43' = '  +3.29%  '
    'D44' = '5.21'
    'E44' = '  +6.68%  '
    'D45' = '276.17'
    'E45' = '  +6.57%  '
    'D46' = '0.0934'
    'E46' = '  +1.55%  '
    'D47' = '0.0505'
    'E47' = '  -0.37%  '
    'D48' = '0.559'
    'E48' = '  +1.84%  '
    'D49' = '0.0218'
    'E49' = '  +1.96%  '
    'D50' = '0.380'
    'E50' = '  +1.60%  '
    'D51' = '17.76'
    'E51' = '  +8.04%  '
}

foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
    $cell.Style = "Normal"
}
